$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-01-15 Monday" "2024-01-16 Tuesday"

Replace-Text "17×57=969" "24×82=1968"
Replace-Text "81×31=2511" "27×11=297"
Replace-Text "13×65=845" "98×65=6370"
Replace-Text "48×94=4512" "65×91=5915"
Replace-Text "23×96=2208" "94×25=2350"

Replace-Text "47×48=2256" "75×11=825"
Replace-Text "50×66=3300" "23×28=644"
Replace-Text "40×95=3800" "38×69=2622"
Replace-Text "77×48=3696" "43×77=3311"
Replace-Text "37×73=2701" "92×68=6256"

Replace-Text "96×98=9408" "94×18=1692"
Replace-Text "28×93=2604" "39×49=1911"
Replace-Text "12×92=1104" "87×46=4002"
Replace-Text "60×68=4080" "36×54=1944"
Replace-Text "81×95=7695" "58×97=5626"

Replace-Text "29×40=1160" "47×50=2350"
Replace-Text "60×69=4140" "59×39=2301"
Replace-Text "14×69=966" "96×48=4608"
Replace-Text "20×19=380" "68×74=5032"
Replace-Text "30×19=570" "43×17=731"

Replace-Text "29×89=2581" "11×78=858"
Replace-Text "14×19=266" "25×45=1125"
Replace-Text "13×93=1209" "40×75=3000"
Replace-Text "11×31=341" "66×78=5148"
Replace-Text "47×44=2068" "61×28=1708"
